$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 60, shifting existing rows 60:134 down to 61:135.
$ws.Rows("60:60").Insert()

# Populate the newly inserted row 60 with the new data record.
$ws.Range("A60").Value = 7
$ws.Range("B60").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C60").Value = "Ñuble"
$ws.Range("D60").Value = 44413
$ws.Range("E60").Value = 16
$ws.Range("F60").Value = 100112023
$ws.Range("G60").Value = "Brócoli"
$ws.Range("H60").Value = "Sin especificar"
$ws.Range("I60").Value = "Primera"
$ws.Range("J60").Value = 300
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 650
$ws.Range("M60").Value = 625
$ws.Range("N60").Value = "$/unidad"
$ws.Range("O60").Value = "Provincia de Diguillín"
$ws.Range("P60").Value = 625
$ws.Range("Q60").Value = 1
$ws.Range("R60").Value = "Hortaliza"
